# This script applies a row-data rotation to the "Damasco" sheet.
# The values originally on rows 2,3,5,6,7,8 (columns D, L-T) get redistributed
# among those same rows (row 4 is left untouched), matching the target diff:
#   new row 2 <= old row 8
#   new row 3 <= old row 7
#   new row 4 <= old row 4 (unchanged)
#   new row 5 <= old row 2
#   new row 6 <= old row 5
#   new row 7 <= old row 3
#   new row 8 <= old row 6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the "before" values for the affected columns on every affected row.
$rows = @(2, 3, 5, 6, 7, 8)
$data = @{}
foreach ($r in $rows) {
    $data[$r] = @{
        D = $ws.Cells.Item($r, 4).Value()
        L = $ws.Cells.Item($r, 12).Value()
        M = $ws.Cells.Item($r, 13).Value()
        N = $ws.Cells.Item($r, 14).Value()
        O = $ws.Cells.Item($r, 15).Value()
        P = $ws.Cells.Item($r, 16).Value()
        Q = $ws.Cells.Item($r, 17).Value()
        R = $ws.Cells.Item($r, 18).Value()
        S = $ws.Cells.Item($r, 19).Value()
        T = $ws.Cells.Item($r, 20).Value()
    }
}

# Mapping of destination row -> source row (data to copy from).
$mapping = @{
    2 = 8
    3 = 7
    5 = 2
    6 = 5
    7 = 3
    8 = 6
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $src = $data[$srcRow]

    $ws.Cells.Item($destRow, 4).Value = $src.D
    $ws.Cells.Item($destRow, 12).Value = $src.L
    $ws.Cells.Item($destRow, 13).Value = $src.M
    $ws.Cells.Item($destRow, 14).Value = $src.N
    $ws.Cells.Item($destRow, 15).Value = $src.O
    $ws.Cells.Item($destRow, 16).Value = $src.P
    $ws.Cells.Item($destRow, 17).Value = $src.Q
    $ws.Cells.Item($destRow, 18).Value = $src.R
    $ws.Cells.Item($destRow, 19).Value = $src.S
    $ws.Cells.Item($destRow, 20).Value = $src.T
}
